$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the row above (D7 -> D8) to get font/border inheritance close to original,
# then overwrite values + styles explicitly to match target structure.

$ws.Range("A8").Value = "0.1uF"
$ws.Range("B8").Value = "C2"
$ws.Range("C8").Value = "C0603"
$ws.Range("D8").Value = "C1591"

$ws.Range("A8").Font.Name = "맑은 고딕"
$ws.Range("A8").HorizontalAlignment = -4131  # xlLeft
$ws.Range("B8").Font.Name = "맑은 고딕"
$ws.Range("C8").Font.Name = "맑은 고딕"

$ws.Range("D8").Font.Name = "Arial"
$ws.Range("D8").Borders.LineStyle = 1
$ws.Range("D8").Borders.Weight = 2
$ws.Range("D8").Borders.Item(8).LineStyle = -4142
$ws.Range("D8").Borders.Item(9).LineStyle = -4142

$ws.Range("D11").Select()
